$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the row for product id 6775928 ("Mezeast Falafel Wrap Kit"); all
#    rows below it (382-394) shift up by one, becoming rows 381-393 and the
#    sheet dimension becomes A1:O393.
$ws.Rows.Item(381).Delete()

# 2) A handful of products toggled their "Online kein Bestand" (out of stock)
#    marker in the aria-label column M. Row numbers below refer to the
#    *post-deletion* layout (these rows are all above the deleted row, so
#    their numbers are unchanged).
$ws.Range("M3").Value = "Prix Garantie Toast - Online kein Bestand 2.15 Schweizer Franken"
$ws.Range("M10").Value = "Ölz Vollkorn Sandwich Toast Soft - Online kein Bestand 4.35 Schweizer Franken"
$ws.Range("M12").Value = "Naturaplan Bio Vollkorntoast 10 Scheiben - Online kein Bestand 2.50 Schweizer Franken"
$ws.Range("M360").Value = "Betty Bossi Dinkel Pizzateig rechteckig 3.30 Schweizer Franken"

# 3) The crawl was re-run later the same day, so every row's timestamp
#    (column O) is refreshed to the new crawl time.
$lastRow = $ws.Range("A1").End(-4121).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 15).Value = "2023-03-06 12:57:43"
}
